# Adds 39 new Customer/Account rows (sheet rows 48-86) to match the
# updated CUSTOMERS export: new TC/Customer_ID/PD rows reusing the existing
# TC and PD codes, plus 39 brand-new Customer_ID values.
#
# Columns A (TC) and C (PD) reuse values that already exist on the sheet, so
# each new cell is populated by Copy + PasteSpecial(Values) from the first
# cell on the sheet holding that value - this keeps the cell typed as a
# shared string (matching the rest of the sheet) without Excel re-inferring
# it as a number and without introducing a new cell style.
#
# Column B (Customer_ID) holds brand-new numeric-looking IDs. Writing them
# straight to .Value would make Excel infer a number, so each is first
# produced as a text formula result in a scratch cell (Z1, well outside the
# used range), copied, and pasted as values into the target cell - this
# lands it as a plain shared string too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 48; ASrc = "A26"; B = "17704760"; CSrc = "C26" },
    @{ Row = 49; ASrc = "A4"; B = "17704761"; CSrc = "C27" },
    @{ Row = 50; ASrc = "A28"; B = "17704762"; CSrc = "C28" },
    @{ Row = 51; ASrc = "A2"; B = "17704763"; CSrc = "C28" },
    @{ Row = 52; ASrc = "A4"; B = "17704764"; CSrc = "C30" },
    @{ Row = 53; ASrc = "A26"; B = "17704767"; CSrc = "C26" },
    @{ Row = 54; ASrc = "A4"; B = "17704768"; CSrc = "C27" },
    @{ Row = 55; ASrc = "A28"; B = "17704769"; CSrc = "C28" },
    @{ Row = 56; ASrc = "A26"; B = "17704772"; CSrc = "C26" },
    @{ Row = 57; ASrc = "A4"; B = "17704773"; CSrc = "C27" },
    @{ Row = 58; ASrc = "A28"; B = "17704774"; CSrc = "C28" },
    @{ Row = 59; ASrc = "A26"; B = "17704775"; CSrc = "C26" },
    @{ Row = 60; ASrc = "A4"; B = "17704776"; CSrc = "C27" },
    @{ Row = 61; ASrc = "A28"; B = "17704777"; CSrc = "C28" },
    @{ Row = 62; ASrc = "A4"; B = "17704779"; CSrc = "C27" },
    @{ Row = 63; ASrc = "A4"; B = "17704781"; CSrc = "C27" },
    @{ Row = 64; ASrc = "A4"; B = "17704783"; CSrc = "C27" },
    @{ Row = 65; ASrc = "A4"; B = "17704786"; CSrc = "C27" },
    @{ Row = 66; ASrc = "A4"; B = "17704787"; CSrc = "C27" },
    @{ Row = 67; ASrc = "A4"; B = "17704789"; CSrc = "C27" },
    @{ Row = 68; ASrc = "A26"; B = "17704791"; CSrc = "C26" },
    @{ Row = 69; ASrc = "A28"; B = "17704793"; CSrc = "C28" },
    @{ Row = 70; ASrc = "A2"; B = "17704794"; CSrc = "C28" },
    @{ Row = 71; ASrc = "A4"; B = "17704795"; CSrc = "C30" },
    @{ Row = 72; ASrc = "A26"; B = "17704800"; CSrc = "C26" },
    @{ Row = 73; ASrc = "A4"; B = "17704801"; CSrc = "C27" },
    @{ Row = 74; ASrc = "A28"; B = "17704802"; CSrc = "C28" },
    @{ Row = 75; ASrc = "A2"; B = "17704803"; CSrc = "C28" },
    @{ Row = 76; ASrc = "A4"; B = "17704804"; CSrc = "C30" },
    @{ Row = 77; ASrc = "A26"; B = "17704808"; CSrc = "C26" },
    @{ Row = 78; ASrc = "A4"; B = "17704809"; CSrc = "C27" },
    @{ Row = 79; ASrc = "A28"; B = "17704810"; CSrc = "C28" },
    @{ Row = 80; ASrc = "A2"; B = "17704811"; CSrc = "C28" },
    @{ Row = 81; ASrc = "A4"; B = "17704812"; CSrc = "C30" },
    @{ Row = 82; ASrc = "A7"; B = "17704813"; CSrc = "C31" },
    @{ Row = 83; ASrc = "A4"; B = "17704815"; CSrc = "C32" },
    @{ Row = 84; ASrc = "A2"; B = "17704816"; CSrc = "C33" },
    @{ Row = 85; ASrc = "A7"; B = "17704817"; CSrc = "C34" },
    @{ Row = 86; ASrc = "A7"; B = "17704818"; CSrc = "C2" }
)

foreach ($row in $newRows) {
    $r = $row.Row

    $ws.Range($row.ASrc).Copy()
    $ws.Range("A$r").PasteSpecial("Values")

    $ws.Range("Z1").Formula = "=""" + $row.B + """"
    $ws.Range("Z1").Copy()
    $ws.Range("B$r").PasteSpecial("Values")

    $ws.Range($row.CSrc).Copy()
    $ws.Range("C$r").PasteSpecial("Values")
}

# Clean up the scratch cell so it leaves no trace in the saved sheet.
$ws.Range("Z1").Value = ""

$ws.Range("A1").Select()
